$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("page33_table1")
Write-Host $ws.Name
$ws.Range("A1").Value = "test"
